$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 0.0002083333333333333
$ws.Range("K2").Value = 1426
$ws.Range("L2").Value = 0.002852
